$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'312.69"
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.Value = "'1.72%"
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.Value = "'39.96"
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.Value = "'-2.72%"
$c.Style = 'Normal'

$c = $ws.Range('D4')
$c.Value = "'5.189"
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.Value = "'-1.00%"
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.Value = "'0.07580"
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.Value = "'-1.01%"
$c.Style = 'Normal'

$c = $ws.Range('B6')
$c.Value = "'GateToken"
$c.Style = 'Normal'
$c = $ws.Range('C6')
$c.Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.Value = "'4.325"
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.Value = "'-0.04%"
$c.Style = 'Normal'

$c = $ws.Range('B7')
$c.Value = "'FTXToken"
$c.Style = 'Normal'
$c = $ws.Range('C7')
$c.Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.Value = "'1.661"
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.Value = "'2.22%"
$c.Style = 'Normal'

$c = $ws.Range('B8')
$c.Value = "'MXToken"
$c.Style = 'Normal'
$c = $ws.Range('C8')
$c.Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.Value = "'0.9259"
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.Value = "'0.97%"
$c.Style = 'Normal'

$c = $ws.Range('B9')
$c.Value = "'BTSEToken"
$c.Style = 'Normal'
$c = $ws.Range('C9')
$c.Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.Value = "'2.424"
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.Value = "'-0.65%"
$c.Style = 'Normal'

$c = $ws.Range('B10')
$c.Value = "'LiechtensteinCryptoassetsExchange"
$c.Style = 'Normal'
$c = $ws.Range('C10')
$c.Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.Value = "'0.1200"
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.Value = "'-4.38%"
$c.Style = 'Normal'

$c = $ws.Range('B11')
$c.Value = "'WazirX"
$c.Style = 'Normal'
$c = $ws.Range('C11')
$c.Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.Value = "'0.1827"
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.Value = "'0.09%"
$c.Style = 'Normal'

$c = $ws.Range('B12')
$c.Value = "'MandalaExchangeToken"
$c.Style = 'Normal'
$c = $ws.Range('C12')
$c.Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.Value = "'0.09046"
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.Value = "'-1.84%"
$c.Style = 'Normal'

$c = $ws.Range('B13')
$c.Value = "'BitrueCoin"
$c.Style = 'Normal'
$c = $ws.Range('C13')
$c.Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.Value = "'0.04168"
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.Value = "'-2.68%"
$c.Style = 'Normal'

$c = $ws.Range('B14')
$c.Value = "'BitMartToken"
$c.Style = 'Normal'
$c = $ws.Range('C14')
$c.Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.Value = "'0.1053"
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.Value = "'0.12%"
$c.Style = 'Normal'

$c = $ws.Range('B15')
$c.Value = "'BitForexToken"
$c.Style = 'Normal'
$c = $ws.Range('C15')
$c.Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.Value = "'0.001296"
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.Value = "'2.22%"
$c.Style = 'Normal'

$c = $ws.Range('B16')
$c.Value = "'TigerCash"
$c.Style = 'Normal'
$c = $ws.Range('C16')
$c.Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.Value = "'0.005839"
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.Value = "'0.02%"
$c.Style = 'Normal'

$c = $ws.Range('B17')
$c.Value = "'UpBots"
$c.Style = 'Normal'
$c = $ws.Range('C17')
$c.Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.Value = "'0.007522"
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.Value = "'0.18%"
$c.Style = 'Normal'

$c = $ws.Range('B18')
$c.Value = "'LEO"
$c.Style = 'Normal'
$c = $ws.Range('C18')
$c.Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.Value = "'3.349"
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.Value = "'-0.20%"
$c.Style = 'Normal'

$c = $ws.Range('E19')
$c.Value = "'0.61%"
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.Value = "'7.588"
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.Value = "'6.04%"
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.Value = "'0.1352"
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.Value = "'-2.48%"
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.Value = "'0.2811"
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.Value = "'0.04018"
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.Value = "'-1.30%"
$c.Style = 'Normal'

$c = $ws.Range('D24')
$c.Value = "'0.001269"
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.Value = "'0.43%"
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.Value = "'0.004069"
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.Value = "'-1.87%"
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.Value = "'-0.25%"
$c.Style = 'Normal'

$c = $ws.Range('D38')
$c.Value = "'0.02421"
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.Value = "'-1.64%"
$c.Style = 'Normal'

$c = $ws.Range('D39')
$c.Value = "'0.05159"
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.Value = "'-2.39%"
$c.Style = 'Normal'

$c = $ws.Range('D40')
$c.Value = "'0.007748"
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.Value = "'-1.35%"
$c.Style = 'Normal'

$c = $ws.Range('D41')
$c.Value = "'0.1299"
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.Value = "'-1.13%"
$c.Style = 'Normal'

$c = $ws.Range('D42')
$c.Value = "'0.007623"
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.Value = "'11.56%"
$c.Style = 'Normal'

$c = $ws.Range('D43')
$c.Value = "'0.003302"
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.Value = "'72.69%"
$c.Style = 'Normal'

$c = $ws.Range('D44')
$c.Value = "'0.008201"
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.Value = "'5.98%"
$c.Style = 'Normal'

$c = $ws.Range('D45')
$c.Value = "'0.3103"
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.Value = "'1.74%"
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.Value = "'0.00006586"
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.Value = "'-2.06%"
$c.Style = 'Normal'

$c = $ws.Range('E47')
$c.Value = "'-0.24%"
$c.Style = 'Normal'

$c = $ws.Range('D48')
$c.Value = "'0.2803"
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.Value = "'36.43%"
$c.Style = 'Normal'

$c = $ws.Range('D49')
$c.Value = "'0.004202"
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.Value = "'2.49%"
$c.Style = 'Normal'

$c = $ws.Range('D50')
$c.Value = "'0.00002101"
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.Value = "'-0.24%"
$c.Style = 'Normal'

$c = $ws.Range('D51')
$c.Value = "'0.0002001"
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.Value = "'-0.24%"
$c.Style = 'Normal'
